$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("D6").Value = "[Python] Object Detection Mosaic Augmentation :: YOLO v5"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Python-Object-Detection-Mosaic-Augmentation-YOLO-v5"

# Row 9
$ws.Range("D9").Value = "수학적으로 그렇게 큰 도전이 아니라는게 좀 많이 알려졌으면 좋겠다"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/student-talks/#utm_source=rss&utm_medium=rss&utm_campaign=student-talks"

# Row 26
$ws.Range("D26").Value = "생성 모델의 새로운 흐름 확산 모델(Diffusion model)에 관하여"

# Row 27
$ws.Range("D27").Value = "Apache Beam으로 머신러닝 데이터 파이프라인 구축하기 2편 - 개발 및 최적화"
$ws.Range("E27").Value = "https://blog.pingpong.us/apache-beam-2/"

# Row 42
$ws.Range("D42").Value = "Python SYS 파라미터 및 변수"
$ws.Range("E42").Value = "https://kjk92.tistory.com/86"

# Row 51
$ws.Range("D51").Value = "1일 1포스팅은 정말 중요한가? 누적 방문수 3백만 블로그 분석 결과 공유"
$ws.Range("E51").Value = "https://bskyvision.com/1285"
